# Update sheet1 ("ランサーズ") with the 2025-10-25 01:15:18 scrape data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks up front; they will be recreated below once
# all cell values (and the new URLs) are in their final positions.
$ws.Hyperlinks.Delete()

$ws.Cells.Item(2,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(2,2).Value = '自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5408668'
$ws.Cells.Item(2,7).Value = 305
$ws.Cells.Item(2,8).Value = '🔥Python ◆開発 ○PHP'

$ws.Cells.Item(3,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(3,2).Value = '【時給1,600円 / 学生限定】AIでプロダクトを生成したことがある学生の方を大募集!!'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5420120'
$ws.Cells.Item(3,7).Value = 303
$ws.Cells.Item(3,8).Value = '🔥AI,Ai'

$ws.Cells.Item(4,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(4,2).Value = '製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5419380'
$ws.Cells.Item(4,7).Value = 298
$ws.Cells.Item(4,8).Value = '🔥AI,Ai'

$ws.Cells.Item(5,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(5,2).Value = '【急募】経験豊富な業務システム開発パートナーを募集'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5419860'
$ws.Cells.Item(5,7).Value = 125
$ws.Cells.Item(5,8).Value = '◆開発,システム開発'

$ws.Cells.Item(6,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(6,2).Value = '【低コスト】住宅リフォーム見積依頼自動化システム構築'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5412955'
$ws.Cells.Item(6,7).Value = 110
$ws.Cells.Item(6,8).Value = '◆自動化'

$ws.Cells.Item(7,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(7,2).Value = '【学生発スタートアップ】留学×住まいマッチングアプリ開発仲間募集'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5420198'
$ws.Cells.Item(7,7).Value = 100
$ws.Cells.Item(7,8).Value = '◆開発 ◇アプリ'

$ws.Cells.Item(8,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(8,2).Value = 'Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5419636'
$ws.Cells.Item(8,7).Value = 85
$ws.Cells.Item(8,8).Value = '★Java'

$ws.Cells.Item(9,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(9,2).Value = 'UIPATHのシステムの開発'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '10,000 円 ~'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5419904'
$ws.Cells.Item(9,7).Value = 75
$ws.Cells.Item(9,8).Value = '◆開発'

$ws.Cells.Item(10,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(10,2).Value = 'IB報酬を得るための高性能EA開発依頼'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5419587'
$ws.Cells.Item(10,7).Value = 68
$ws.Cells.Item(10,8).Value = '◆開発'

$ws.Cells.Item(11,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(11,2).Value = 'クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5419638'
$ws.Cells.Item(11,7).Value = 38
$ws.Cells.Item(11,8).Value = '◇管理'

$ws.Cells.Item(12,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(12,2).Value = '【Ubuntu】MySQLデータを自動CSV化しクラウド保存構築'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5420180'
$ws.Cells.Item(12,7).Value = 30
$ws.Cells.Item(12,8).Value = '◇MySQL'

$ws.Cells.Item(13,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(13,2).Value = '【緊急】ロリポップ Wordpress リダイレクトハッキング復旧依頼'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '20,000 円 ~ 30,000 円 / 募集期間 5 日、取引期間 0 日'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5419656'
$ws.Cells.Item(13,7).Value = 25
$ws.Cells.Item(13,8).Value = '○WordPress'

$ws.Cells.Item(14,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(14,2).Value = '【Braze経験者募集】CRM/マーケティングオートメーション支援(中級者以上)'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5419838'
$ws.Cells.Item(14,7).Value = 25
$ws.Cells.Item(14,8).Value = ""

$ws.Cells.Item(15,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(15,2).Value = '【KARTE経験者募集】CX改善/Web接客施策の設計・実装(中級者以上)'
$ws.Cells.Item(15,3).Value = 'システム開発'
$ws.Cells.Item(15,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(15,5).Value = '期限情報なし'
$ws.Cells.Item(15,6).Value = 'https://www.lancers.jp/work/detail/5419829'
$ws.Cells.Item(15,7).Value = 25
$ws.Cells.Item(15,8).Value = ""

$ws.Cells.Item(16,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(16,2).Value = '注目 【急募】YouTubeの音楽配信構築の依頼です'
$ws.Cells.Item(16,3).Value = 'システム開発'
$ws.Cells.Item(16,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(16,5).Value = '期限情報なし'
$ws.Cells.Item(16,6).Value = 'https://www.lancers.jp/work/detail/5420233'
$ws.Cells.Item(16,7).Value = 13
$ws.Cells.Item(16,8).Value = ""

$ws.Cells.Item(17,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(17,2).Value = '【急募】イベント用問い合わせLINE構築のフリーランス募集!'
$ws.Cells.Item(17,3).Value = 'システム開発'
$ws.Cells.Item(17,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(17,5).Value = '期限情報なし'
$ws.Cells.Item(17,6).Value = 'https://www.lancers.jp/work/detail/5420186'
$ws.Cells.Item(17,7).Value = 10
$ws.Cells.Item(17,8).Value = ""

$ws.Cells.Item(18,1).Value = '2025-10-25 01:15:18'
$ws.Cells.Item(18,2).Value = '【急募】Google Play Consoleでのクローズテスト実施者募集!'
$ws.Cells.Item(18,3).Value = 'システム開発'
$ws.Cells.Item(18,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(18,5).Value = '期限情報なし'
$ws.Cells.Item(18,6).Value = 'https://www.lancers.jp/work/detail/5419425'
$ws.Cells.Item(18,7).Value = 10
$ws.Cells.Item(18,8).Value = ""

# Re-create one hyperlink per row (F2:F18), pointing at the URL now shown
# in that cell, then restore the shared "Hyperlink" cell style on column F.
for ($r = 2; $r -le 18; $r++) {
    $urlCell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($urlCell, $urlCell.Value) | Out-Null
}
$ws.Range("F2:F18").Style = "Hyperlink"

$ws.Range("A1").Select() | Out-Null